$wb = $excel.ActiveWorkbook

# ---- Sheet1: updated computed/source values ----
$s1 = $wb.Worksheets.Item("Sheet1")
$s1.Range("A1").Value = "5,3572"
$s1.Range("A2").Value = "4,4004"
$s1.Range("A3").Value = "4,6363"
$s1.Range("A4").Value = "1874,8"
$s1.Range("A7").Value = "22,5"
$s1.Range("A8").Value = "26438,67"
# A9 would otherwise be auto-parsed as the number 687 (thousands-separator
# grouping matches "000" + "687"); keep it as literal text like the source.
$s1.Range("A9").Value = "'0,000687"
$s1.Range("A10").Value = "1567,04"
$s1.Range("A12").Value = "0,1421"
$s1.Range("A13").Value = "143,76"
$s1.Range("A34").Value = "6662,0"

# ---- data sheet: relabel/clear entries whose lookups now error out ----
$data = $wb.Worksheets.Item("data")

# Row 1: tickers / labels
$data.Range("A1").Value = "swda-etf"
$data.Range("C1").Value = "-"
$data.Range("G1").Value = "-"
$data.Range("I1").Value = "-"
$data.Range("J1").Value = "-"
$data.Range("K1").Value = "-"
$data.Range("L1").Value = "-"
$data.Range("M1").Value = "-"

# Row 2: sheet-name references - clear the ones that now error
$data.Range("C2").Value = ""
$data.Range("G2").Value = ""
$data.Range("I2").Value = ""
$data.Range("J2").Value = ""
$data.Range("K2").Value = ""
$data.Range("L2").Value = ""
$data.Range("M2").Value = ""

# Row 3: cell references - A3 now points at A34 (instead of A1); others cleared
$data.Range("A3").Value = "A34"
$data.Range("C3").Value = ""
$data.Range("G3").Value = ""
$data.Range("I3").Value = ""
$data.Range("J3").Value = ""
$data.Range("K3").Value = ""
$data.Range("L3").Value = ""
$data.Range("M3").Value = ""

# Row 4: currency units - a GBP column is inserted before the old A4,
# shifting old values right by one and changing what is now C4 to USD
$data.Range("A4").Value = "GBP"
$data.Range("B4").Value = "PLN"
$data.Range("C4").Value = "USD"
